$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.565.33'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '2.393.65'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.535'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.34%  '
$ws.Range('D9').Value = '2.398.50'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.108'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('E11').Value = '  -0.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.18'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').Value = '2.834.54'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('E16').Value = '  -1.73%  '
$ws.Range('D17').Value = '60.231.46'
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').Value = '2.413.30'
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.15'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +10.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  -2.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '64.79'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '568.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.43%  '
$ws.Range('D29').Value = '2.504.25'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').Value = '0.0₃0935'
$ws.Range('E30').Value = '  +1.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.07'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.33'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.82'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('E34').Value = '  -2.00%  '
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('E36').Value = '  +3.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '152.12'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('E38').Value = '  +0.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.59'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.27'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.14'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.50'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.62'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.67'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('D46').Value = '0.0₆0286'
$ws.Range('E46').Value = '  +3.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '141.86'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.56'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.589'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.32'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.39%  '
